$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 44.29505033333334
$ws.Range("H2").Value = 132.885151
$ws.Range("I2").Value = 0.9830698162761968
$ws.Range("J2").Value = 0.9830698162761969
$ws.Range("M2").Value = 7.579746333333333
$ws.Range("N2").Value = 22.739239
$ws.Range("O2").Value = 0.0686314777863378
$ws.Range("P2").Value = 0.0686314777863378
$ws.Range("Q2").Value = 335.7452453488988
$ws.Range("R2").Value = 3021.707208140089
$ws.Range("S2").Value = 0.06746953425817898
$ws.Range("T2").Value = 0.06746953425817899

$ws.Range("G3").Value = 44.29505033333334
$ws.Range("H3").Value = 132.885151
$ws.Range("I3").Value = 0.9830698162761968
$ws.Range("J3").Value = 0.9830698162761969
$ws.Range("O3").Value = 0.0596740760116217
$ws.Range("P3").Value = 0.05967407601162171
$ws.Range("Q3").Value = 291.9256285558099
$ws.Range("R3").Value = 2627.330657002289
$ws.Range("S3").Value = 0.05866378294119675
$ws.Range("T3").Value = 0.05866378294119676

$ws.Range("G4").Value = 44.29505033333334
$ws.Range("H4").Value = 132.885151
$ws.Range("I4").Value = 0.9830698162761968
$ws.Range("J4").Value = 0.9830698162761969
$ws.Range("M4").Value = 96.08192699999999
$ws.Range("N4").Value = 288.245781
$ws.Range("O4").Value = 0.8699822327258658
$ws.Range("P4").Value = 0.8699822327258659
$ws.Range("Q4").Value = 4255.953792588659
$ws.Range("R4").Value = 38303.58413329793
$ws.Range("S4").Value = 0.8552532736893724
$ws.Range("T4").Value = 0.8552532736893727

$ws.Range("G5").Value = 44.29505033333334
$ws.Range("H5").Value = 132.885151
$ws.Range("I5").Value = 0.9830698162761968
$ws.Range("J5").Value = 0.9830698162761969
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.189099
$ws.Range("N5").Value = 0.5672970000000001
$ws.Range("O5").Value = 0.001712213476174646
$ws.Range("P5").Value = 0.001712213476174646
$ws.Range("Q5").Value = 8.376149722983001
$ws.Range("R5").Value = 75.38534750684701
$ws.Range("S5").Value = 0.001683225387448637
$ws.Range("T5").Value = 0.001683225387448637

$ws.Range("I6").Value = 0.006814145293655052
$ws.Range("J6").Value = 0.006814145293655053
$ws.Range("M6").Value = 7.579746333333333
$ws.Range("N6").Value = 22.739239
$ws.Range("O6").Value = 0.0686314777863378
$ws.Range("P6").Value = 0.0686314777863378
$ws.Range("Q6").Value = 2.327217096469667
$ws.Range("R6").Value = 20.944953868227
$ws.Range("S6").Value = 0.000467664861354365
$ws.Range("T6").Value = 0.000467664861354365

$ws.Range("I7").Value = 0.006814145293655052
$ws.Range("J7").Value = 0.006814145293655053
$ws.Range("O7").Value = 0.0596740760116217
$ws.Range("P7").Value = 0.05967407601162171
$ws.Range("S7").Value = 0.0004066278242078059
$ws.Range("T7").Value = 0.000406627824207806

$ws.Range("I8").Value = 0.006814145293655052
$ws.Range("J8").Value = 0.006814145293655053
$ws.Range("M8").Value = 96.08192699999999
$ws.Range("N8").Value = 288.245781
$ws.Range("O8").Value = 0.8699822327258658
$ws.Range("P8").Value = 0.8699822327258659
$ws.Range("Q8").Value = 29.500130128737
$ws.Range("R8").Value = 265.501171158633
$ws.Range("S8").Value = 0.005928185336692473
$ws.Range("T8").Value = 0.005928185336692475

$ws.Range("I9").Value = 0.006814145293655052
$ws.Range("J9").Value = 0.006814145293655053
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.189099
$ws.Range("N9").Value = 0.5672970000000001
$ws.Range("O9").Value = 0.001712213476174646
$ws.Range("P9").Value = 0.001712213476174646
$ws.Range("Q9").Value = 0.058059255069
$ws.Range("R9").Value = 0.522533295621
$ws.Range("S9").Value = 0.00001166727140040822
$ws.Range("T9").Value = 0.00001166727140040822

$ws.Range("G10").Value = 0.3685326666666667
$ws.Range("H10").Value = 1.105598
$ws.Range("I10").Value = 0.008179093108268589
$ws.Range("J10").Value = 0.008179093108268589
$ws.Range("M10").Value = 7.579746333333333
$ws.Range("N10").Value = 22.739239
$ws.Range("O10").Value = 0.0686314777863378
$ws.Range("P10").Value = 0.0686314777863378
$ws.Range("Q10").Value = 2.793384128880222
$ws.Range("R10").Value = 25.140457159922
$ws.Range("S10").Value = 0.0005613432469725242
$ws.Range("T10").Value = 0.0005613432469725242

$ws.Range("G11").Value = 0.3685326666666667
$ws.Range("H11").Value = 1.105598
$ws.Range("I11").Value = 0.008179093108268589
$ws.Range("J11").Value = 0.008179093108268589
$ws.Range("O11").Value = 0.0596740760116217
$ws.Range("P11").Value = 0.05967407601162171
$ws.Range("Q11").Value = 2.428807046169111
$ws.Range("R11").Value = 21.859263415522
$ws.Range("S11").Value = 0.000488079823848951
$ws.Range("T11").Value = 0.000488079823848951

$ws.Range("G12").Value = 0.3685326666666667
$ws.Range("H12").Value = 1.105598
$ws.Range("I12").Value = 0.008179093108268589
$ws.Range("J12").Value = 0.008179093108268589
$ws.Range("M12").Value = 96.08192699999999
$ws.Range("N12").Value = 288.245781
$ws.Range("O12").Value = 0.8699822327258658
$ws.Range("P12").Value = 0.8699822327258659
$ws.Range("Q12").Value = 35.409328775782
$ws.Range("R12").Value = 318.683958982038
$ws.Range("S12").Value = 0.007115665684004248
$ws.Range("T12").Value = 0.007115665684004249

$ws.Range("G13").Value = 0.3685326666666667
$ws.Range("H13").Value = 1.105598
$ws.Range("I13").Value = 0.008179093108268589
$ws.Range("J13").Value = 0.008179093108268589
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.189099
$ws.Range("N13").Value = 0.5672970000000001
$ws.Range("O13").Value = 0.001712213476174646
$ws.Range("P13").Value = 0.001712213476174646
$ws.Range("Q13").Value = 0.06968915873400001
$ws.Range("R13").Value = 0.6272024286060001
$ws.Range("S13").Value = 0.00001400435344286465
$ws.Range("T13").Value = 0.00001400435344286465

$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.08727466666666667
$ws.Range("H14").Value = 0.261824
$ws.Range("I14").Value = 0.001936945321879485
$ws.Range("J14").Value = 0.001936945321879485
$ws.Range("M14").Value = 7.579746333333333
$ws.Range("N14").Value = 22.739239
$ws.Range("O14").Value = 0.0686314777863378
$ws.Range("P14").Value = 0.0686314777863378
$ws.Range("Q14").Value = 0.6615198346595556
$ws.Range("R14").Value = 5.953678511936
$ws.Range("S14").Value = 0.0001329354198319228
$ws.Range("T14").Value = 0.0001329354198319228

$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.08727466666666667
$ws.Range("H15").Value = 0.261824
$ws.Range("I15").Value = 0.001936945321879485
$ws.Range("J15").Value = 0.001936945321879485
$ws.Range("O15").Value = 0.0596740760116217
$ws.Range("P15").Value = 0.05967407601162171
$ws.Range("Q15").Value = 0.5751819160817778
$ws.Range("R15").Value = 5.176637244736
$ws.Range("S15").Value = 0.0001155854223681914
$ws.Range("T15").Value = 0.0001155854223681915

$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.08727466666666667
$ws.Range("H16").Value = 0.261824
$ws.Range("I16").Value = 0.001936945321879485
$ws.Range("J16").Value = 0.001936945321879485
$ws.Range("M16").Value = 96.08192699999999
$ws.Range("N16").Value = 288.245781
$ws.Range("O16").Value = 0.8699822327258658
$ws.Range("P16").Value = 0.8699822327258659
$ws.Range("Q16").Value = 8.385518151615999
$ws.Range("R16").Value = 75.46966336454399
$ws.Range("S16").Value = 0.001685108015796635
$ws.Range("T16").Value = 0.001685108015796636

$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.08727466666666667
$ws.Range("H17").Value = 0.261824
$ws.Range("I17").Value = 0.001936945321879485
$ws.Range("J17").Value = 0.001936945321879485
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.189099
$ws.Range("N17").Value = 0.5672970000000001
$ws.Range("O17").Value = 0.001712213476174646
$ws.Range("P17").Value = 0.001712213476174646
$ws.Range("Q17").Value = 0.016503552192
$ws.Range("R17").Value = 0.148531969728
$ws.Range("S17").Value = 0.000003316463882735491
$ws.Range("T17").Value = 0.000003316463882735491

Write-Output "applied changes"